# Update cryptos list data (prices and 1h volume changes) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.800.41"
$ws.Range("E2").Value = "  +2.89%  "
$ws.Range("D3").Value = "3.731.97"
$ws.Range("E3").Value = "  +6.78%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.44"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.87"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "3.723.51"
$ws.Range("E7").Value = "  +6.75%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +14.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000409"
$ws.Range("E12").Value = "  +58.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.90"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.56"
$ws.Range("E14").Value = "  +7.66%  "
$ws.Range("D15").Value = "4.298.73"
$ws.Range("E15").Value = "  +5.96%  "
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("D18").Value = "3.717.70"
$ws.Range("E18").Value = "  +6.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  +5.56%  "
$ws.Range("E20").Value = "  +4.80%  "
$ws.Range("D21").Value = "66.755.01"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.93"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.47"
$ws.Range("E23").Value = "  +24.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.80"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.67"
$ws.Range("E26").Value = "  +14.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.24"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  +4.12%  "
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("E31").Value = "  +9.75%  "
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.91"
$ws.Range("E35").Value = "  +6.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.15"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0495"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "0.0₃0760"
$ws.Range("E39").Value = "  +8.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  +30.46%  "
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "28.86"
$ws.Range("E42").Value = "  +34.20%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("E45").Value = "  +32.19%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.13"
$ws.Range("E46").Value = "  +6.36%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "146.97"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.39"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.307"
$ws.Range("E51").Value = "  -1.70%  "
